$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.598.16'
$ws.Range('E2').Value = '  +3.02%  '
$ws.Range('D3').Value = '3.395.63'
$ws.Range('E3').Value = '  +2.41%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '588.97'
$ws.Range('E5').Value = '  +2.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.59'
$ws.Range('E6').Value = '  +4.53%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +1.37%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.197'
$ws.Range('E9').Value = '  +10.71%  '
$ws.Range('E10').Value = '  +2.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '48.90'
$ws.Range('E11').Value = '  +7.03%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000283'
$ws.Range('E12').Value = '  +5.21%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '693.17'
$ws.Range('E13').Value = '  -2.36%  '
$ws.Range('D14').Value = '3.946.00'
$ws.Range('E14').Value = '  +2.22%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.60'
$ws.Range('E15').Value = '  +2.69%  '
$ws.Range('D16').Value = '69.579.31'
$ws.Range('E16').Value = '  +3.00%  '
$ws.Range('D17').Value = '3.394.54'
$ws.Range('E17').Value = '  +1.91%  '
$ws.Range('E18').Value = '  +1.89%  '
$ws.Range('E19').Value = '  +2.21%  '
$ws.Range('E20').Value = '  +4.63%  '
$ws.Range('E21').Value = '  +2.04%  '
$ws.Range('E22').Value = '  +1.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.14'
$ws.Range('E23').Value = '  +1.85%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '104.60'
$ws.Range('E24').Value = '  +6.44%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.98'
$ws.Range('E25').Value = '  +3.30%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.73'
$ws.Range('E26').Value = '  +2.42%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.68'
$ws.Range('E27').Value = '  +3.93%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '34.47'
$ws.Range('E28').Value = '  +3.34%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.71'
$ws.Range('E29').Value = '  +3.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.09'
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '11.20'
$ws.Range('E31').Value = '  +2.43%  '
$ws.Range('B32').Value = 'dogwifhat'
$ws.Range('C32').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.69'
$ws.Range('E32').Value = '  +11.43%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '558.97'
$ws.Range('E33').Value = '  -1.53%  '
$ws.Range('E34').Value = '  +2.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '58.55'
$ws.Range('E35').Value = '  +2.14%  '
$ws.Range('D36').Value = '3.732.93'
$ws.Range('E36').Value = '  +1.08%  '
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('E38').Value = '  +9.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.03'
$ws.Range('E39').Value = '  +2.33%  '
$ws.Range('D40').Value = '0.0₃0713'
$ws.Range('E40').Value = '  +7.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.24'
$ws.Range('E41').Value = '  +2.66%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.68'
$ws.Range('E42').Value = '  +2.76%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.341'
$ws.Range('E43').Value = '  +2.82%  '
$ws.Range('E44').Value = '  +3.98%  '
$ws.Range('E45').Value = '  -1.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.67'
$ws.Range('E46').Value = '  -0.45%  '
$ws.Range('E47').Value = '  +1.86%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.41'
$ws.Range('E48').Value = '  +8.20%  '
$ws.Range('E49').Value = '  -0.25%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.67'
$ws.Range('E50').Value = '  +3.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.63'
$ws.Range('E51').Value = '  -1.65%  '
